$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 99, pushing existing rows 99:116 down to 100:117
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly price record
$ws.Cells.Item(99, 1).Value = 11
$ws.Cells.Item(99, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(99, 3).Value = "Bíobío"
$ws.Cells.Item(99, 4).Value = 44551
$ws.Cells.Item(99, 5).Value = 8
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100103
$ws.Cells.Item(99, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(99, 9).Value = 100103004
$ws.Cells.Item(99, 10).Value = "Durazno"
$ws.Cells.Item(99, 11).Value = "Royal Glory"
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 270
$ws.Cells.Item(99, 14).Value = 12000
$ws.Cells.Item(99, 15).Value = 13000
$ws.Cells.Item(99, 16).Value = 12444
$ws.Cells.Item(99, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(99, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(99, 19).Value = 830
$ws.Cells.Item(99, 20).Value = 15
